$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 42.166668
$ws.Range("I11").Value = 42.166668
$ws.Range("K11").Value = 42.166668
$ws.Range("M11").Value = 97.833332
$ws.Range("H15").Value = 3073.0193
$ws.Range("I15").Value = 3073.0193
$ws.Range("K15").Value = 9219.0579
$ws.Range("M15").Value = -9050.0579
$ws.Range("H53").Value = 496.91666
$ws.Range("I53").Value = 496.45456
$ws.Range("J53").Value = 502
$ws.Range("K53").Value = 496.45456
$ws.Range("L53").Value = 502
$ws.Range("M53").Value = 140.54544
$ws.Range("N53").Value = -1776
$ws.Range("H64").Value = 3041.6287
$ws.Range("I64").Value = 3134.6667
$ws.Range("J64").Value = 3022.3794
$ws.Range("K64").Value = 3134.6667
$ws.Range("L64").Value = 3022.3794
$ws.Range("M64").Value = -2886.6667
$ws.Range("N64").Value = -3518.3794
$ws.Range("H67").Value = 3041.6287
$ws.Range("I67").Value = 3134.6667
$ws.Range("J67").Value = 3022.3794
$ws.Range("K67").Value = 3134.6667
$ws.Range("L67").Value = 3022.3794
$ws.Range("M67").Value = -2276.6667
$ws.Range("N67").Value = -4738.3794
$ws.Range("H80").Value = 1473.069
$ws.Range("I80").Value = 671.375
$ws.Range("J80").Value = 2459.7693
$ws.Range("K80").Value = 2014.125
$ws.Range("L80").Value = 7379.3079
$ws.Range("M80").Value = -1016.125
$ws.Range("N80").Value = -9375.3079
$ws.Range("H83").Value = 1473.069
$ws.Range("I83").Value = 671.375
$ws.Range("J83").Value = 2459.7693
$ws.Range("K83").Value = 6042.375
$ws.Range("L83").Value = 22137.9237
$ws.Range("M83").Value = -1050.375
$ws.Range("N83").Value = -32121.9237
$ws.Range("H137").Value = 1832.2916
$ws.Range("I137").Value = 1650.4
$ws.Range("J137").Value = 2322
$ws.Range("K137").Value = 4951.200000000001
$ws.Range("L137").Value = 6966
$ws.Range("M137").Value = -2401.200000000001
$ws.Range("N137").Value = -12066

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 2003398
$ws.Range("I13").Value = 10000000
$ws.Range("J13").Value = 4247.5
$ws.Range("K13").Value = 10000000
$ws.Range("L13").Value = 4247.5
$ws.Range("M13").Value = -9999856
$ws.Range("N13").Value = -4535.5
$ws.Range("H102").Value = 2445.0952
$ws.Range("I102").Value = 2355.75
$ws.Range("J102").Value = 2564.2222
$ws.Range("K102").Value = 2355.75
$ws.Range("L102").Value = 2564.2222
$ws.Range("M102").Value = -733.75
$ws.Range("N102").Value = -5808.2222

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8335775
$ws.Range("I31").Value = 11113124
$ws.Range("J31").Value = 6669365.5
$ws.Range("K31").Value = 11113124
$ws.Range("L31").Value = 6669365.5
$ws.Range("M31").Value = -11112829
$ws.Range("N31").Value = -6669955.5
$ws.Range("H34").Value = 8335775
$ws.Range("I34").Value = 11113124
$ws.Range("J34").Value = 6669365.5
$ws.Range("K34").Value = 11113124
$ws.Range("L34").Value = 6669365.5
$ws.Range("M34").Value = -11112922
$ws.Range("N34").Value = -6669769.5
$ws.Range("H62").Value = 76925830
$ws.Range("I62").Value = 2874.25
$ws.Range("K62").Value = 2874.25
$ws.Range("M62").Value = -2250.25
$ws.Range("H65").Value = 76925830
$ws.Range("I65").Value = 2874.25
$ws.Range("K65").Value = 14371.25
$ws.Range("M65").Value = -11251.25
$ws.Range("H134").Value = 2837.3428
$ws.Range("I134").Value = 2682.5518
$ws.Range("J134").Value = 3585.5
$ws.Range("K134").Value = 8047.655400000001
$ws.Range("L134").Value = 10756.5
$ws.Range("M134").Value = -5512.655400000001
$ws.Range("N134").Value = -15826.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 715.3461
$ws.Range("I5").Value = 735.96
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 2207.88
$ws.Range("L5").Value = 600
$ws.Range("M5").Value = -2095.88
$ws.Range("N5").Value = -824
$ws.Range("H135").Value = 715.3461
$ws.Range("I135").Value = 735.96
$ws.Range("J135").Value = 200
$ws.Range("K135").Value = 6623.64
$ws.Range("L135").Value = 1800
$ws.Range("M135").Value = -4088.64
$ws.Range("N135").Value = -6870

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 2865.3333
$ws.Range("I21").Value = 2798
$ws.Range("K21").Value = 2798
$ws.Range("M21").Value = -2625
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H30").Value = 2865.3333
$ws.Range("I30").Value = 2798
$ws.Range("K30").Value = 2798
$ws.Range("M30").Value = -2693
$ws.Range("H33").Value = 7500
$ws.Range("J33").Value = 7500
$ws.Range("L33").Value = 7500
$ws.Range("N33").Value = -8004
$ws.Range("H52").Value = 5000
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H55").Value = 2830
$ws.Range("I55").Value = 2830
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 2830
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -2503
$ws.Range("N55").ClearContents()
$ws.Range("H57").Value = 10978.2
$ws.Range("I57").Value = 6000
$ws.Range("J57").Value = 12222.75
$ws.Range("K57").Value = 6000
$ws.Range("L57").Value = 12222.75
$ws.Range("M57").Value = -5180
$ws.Range("N57").Value = -13862.75
$ws.Range("H58").Value = 7750
$ws.Range("I58").Value = 8666.666999999999
$ws.Range("J58").Value = 5000
$ws.Range("K58").Value = 8666.666999999999
$ws.Range("L58").Value = 5000
$ws.Range("M58").Value = -8389.666999999999
$ws.Range("N58").Value = -5554
$ws.Range("H107").Value = 5952938.5
$ws.Range("I107").Value = 11111751
$ws.Range("J107").Value = 462.07693
$ws.Range("K107").Value = 11111751
$ws.Range("L107").Value = 462.07693
$ws.Range("M107").Value = -11109831
$ws.Range("N107").Value = -4302.07693

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1231574.1
$ws.Range("I132").Value = 2240908.5
$ws.Range("J132").Value = 2819.2173
$ws.Range("K132").Value = 6722725.5
$ws.Range("L132").Value = 8457.651899999999
$ws.Range("M132").Value = -6720195.5
$ws.Range("N132").Value = -13517.6519

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1489.2927
$ws.Range("I132").Value = 1223.5625
$ws.Range("J132").Value = 2434.111
$ws.Range("K132").Value = 3670.6875
$ws.Range("L132").Value = 7302.333
$ws.Range("M132").Value = -1140.6875
$ws.Range("N132").Value = -12362.333
